$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2000
$ws.Range("I88").Value = 2000
$ws.Range("K88").Value = 2000
$ws.Range("M88").Value = -1594
$ws.Range("H91").Value = 2000
$ws.Range("I91").Value = 2000
$ws.Range("K91").Value = 2000
$ws.Range("M91").Value = -596
$ws.Range("H98").Value = 550
$ws.Range("I98").Value = 550
$ws.Range("K98").Value = 550
$ws.Range("M98").Value = 948
$ws.Range("H100").Value = 1710.7778
$ws.Range("I100").Value = 1800.125
$ws.Range("K100").Value = 1800.125
$ws.Range("M100").Value = -1259.125
$ws.Range("H122").Value = 550
$ws.Range("I122").Value = 550
$ws.Range("K122").Value = 1650
$ws.Range("M122").Value = 800
$ws.Range("H131").Value = 4472.5
$ws.Range("J131").Value = 4470
$ws.Range("L131").Value = 13410
$ws.Range("N131").Value = -23490
$ws.Range("H137").Value = 2154.7
$ws.Range("I137").Value = 729.4
$ws.Range("K137").Value = 2188.2
$ws.Range("M137").Value = 361.8000000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3631.111
$ws.Range("I45").Value = 2788.6
$ws.Range("K45").Value = 2788.6
$ws.Range("M45").Value = -2411.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 263
$ws.Range("I22").Value = 283.33334
$ws.Range("J22").Value = 202
$ws.Range("K22").Value = 283.33334
$ws.Range("L22").Value = 202
$ws.Range("M22").Value = -110.33334
$ws.Range("N22").Value = -548
$ws.Range("H86").Value = 5899.154
$ws.Range("I86").Value = 3798.8
$ws.Range("J86").Value = 7211.875
$ws.Range("K86").Value = 3798.8
$ws.Range("L86").Value = 7211.875
$ws.Range("M86").Value = -2675.8
$ws.Range("N86").Value = -9457.875
$ws.Range("H89").Value = 5899.154
$ws.Range("I89").Value = 3798.8
$ws.Range("J89").Value = 7211.875
$ws.Range("K89").Value = 18994
$ws.Range("L89").Value = 36059.375
$ws.Range("M89").Value = -13378
$ws.Range("N89").Value = -47291.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2290
$ws.Range("I99").Value = 2148
$ws.Range("K99").Value = 2148
$ws.Range("M99").Value = -650
$ws.Range("H107").Value = 1305.8572
$ws.Range("I107").Value = 497.5
$ws.Range("K107").Value = 497.5
$ws.Range("M107").Value = 1422.5
$ws.Range("H126").Value = 2290
$ws.Range("I126").Value = 2148
$ws.Range("K126").Value = 6444
$ws.Range("M126").Value = -3974
$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 16000.667
$ws.Range("I9").Value = 500
$ws.Range("J9").Value = 23751
$ws.Range("K9").Value = 1500
$ws.Range("L9").Value = 71253
$ws.Range("M9").Value = -1276
$ws.Range("N9").Value = -71701
$ws.Range("H60").Value = 1742.2222
$ws.Range("I60").Value = 185
$ws.Range("J60").Value = 2520.8333
$ws.Range("K60").Value = 555
$ws.Range("L60").Value = 7562.499899999999
$ws.Range("M60").Value = -304
$ws.Range("N60").Value = -8064.499899999999
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H121").Value = 215
$ws.Range("I121").Value = 215
$ws.Range("K121").Value = 645
$ws.Range("M121").Value = 665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 259.5
$ws.Range("I2").Value = 223
$ws.Range("J2").Value = 305.125
$ws.Range("K2").Value = 223
$ws.Range("L2").Value = 305.125
$ws.Range("M2").Value = -110
$ws.Range("N2").Value = -531.125
$ws.Range("H102").Value = 2927.4546
$ws.Range("I102").Value = 2211.5
$ws.Range("J102").Value = 4836.6665
$ws.Range("K102").Value = 2211.5
$ws.Range("L102").Value = 4836.6665
$ws.Range("M102").Value = -589.5
$ws.Range("N102").Value = -8080.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5947.1665
$ws.Range("I7").Value = 5137.8
$ws.Range("J7").Value = 9994
$ws.Range("K7").Value = 5137.8
$ws.Range("L7").Value = 9994
$ws.Range("M7").Value = -5025.8
$ws.Range("N7").Value = -10218
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H97").Value = 35000
$ws.Range("J97").Value = 35000
$ws.Range("L97").Value = 35000
$ws.Range("N97").Value = -36982
$ws.Range("H126").Value = 5947.1665
$ws.Range("I126").Value = 5137.8
$ws.Range("J126").Value = 9994
$ws.Range("K126").Value = 15413.4
$ws.Range("L126").Value = 29982
$ws.Range("M126").Value = -12943.4
$ws.Range("N126").Value = -34922
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H128").Value = 40000
$ws.Range("J128").Value = 40000
$ws.Range("L128").Value = 40000
$ws.Range("N128").Value = -49960
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H131").Value = 250000
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H136").Value = 1297
$ws.Range("I136").Value = 1273.75
$ws.Range("K136").Value = 3821.25
$ws.Range("M136").Value = -1271.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 27000
$ws.Range("J110").Value = 27000
$ws.Range("L110").Value = 27000
$ws.Range("N110").Value = -35180
$ws.Range("H126").Value = 3597.5
$ws.Range("I126").Value = 2257.6086
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 6772.825800000001
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -4302.825800000001
$ws.Range("N126").Value = -28940
$ws.Range("H132").Value = 2209.3333
$ws.Range("I132").Value = 2312
$ws.Range("J132").Value = 1850
$ws.Range("K132").Value = 6936
$ws.Range("L132").Value = 5550
$ws.Range("M132").Value = -4406
$ws.Range("N132").Value = -10610
